# Insert a new data row at row 76 (pushing the existing rows 76..196 down
# to 77..197) and populate it with a new weekly price observation.
#
# The new row reuses the same Mercado/Categoria/etc. values that the sheet
# already uses on every data row, the same Volumen/Precio values that used
# to sit at row 76 (J=160, K/L/M=1500, P=500), and a new Fecha (44495).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 76; this shifts old rows 76-196
# down to 77-197 and copies formatting (incl. the date number format) from
# the row that used to be at 76 (now 77).
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new weekly observation.
$ws.Cells.Item(76, 1).Value = 3
$ws.Cells.Item(76, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(76, 3).Value = "Coquimbo"
$ws.Cells.Item(76, 4).Value = 44495
$ws.Cells.Item(76, 5).Value = 5
$ws.Cells.Item(76, 6).Value = 100112039
$ws.Cells.Item(76, 7).Value = "Ciboulette"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 160
$ws.Cells.Item(76, 11).Value = 1500
$ws.Cells.Item(76, 12).Value = 1500
$ws.Cells.Item(76, 13).Value = 1500
$ws.Cells.Item(76, 14).Value = "$/docena de atados"
$ws.Cells.Item(76, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(76, 16).Value = 500
$ws.Cells.Item(76, 17).Value = 3
$ws.Cells.Item(76, 18).Value = "Hortaliza"
